$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# "Metadata" sheet (sheet1): the StructureDefinition metadata key/value table.
# ---------------------------------------------------------------------------
$meta = $wb.Worksheets.Item("Metadata")

# The source table had a duplicated "Contact" / "No display for ContactDetail"
# row (rows 10 and 11 both held the same text). The new export drops the
# duplicate and instead reports the real Publisher/Jurisdiction metadata, so
# remove the stray duplicate row first (this shifts everything below it up
# by one, matching the new A1:B20 dimension instead of A1:B21).
$meta.Rows.Item(11).Delete()

# Version bump.
$meta.Range("B3").Value = "6.0.0"

# New publication date.
$meta.Range("B8").Value = "2022-01-21T20:46:54+00:00"

# Publisher now has a real value instead of being blank.
$meta.Range("B9").Value = "Alvearie Team"

# Row 10 (formerly the duplicate "Contact" row) now reports Jurisdiction.
$meta.Range("A10").Value = "Jurisdiction"
$meta.Range("B10").Value = "United States of America"

# ---------------------------------------------------------------------------
# "Elements" sheet (sheet2): the element definitions / mappings table.
# ---------------------------------------------------------------------------
$elements = $wb.Worksheets.Item("Elements")

# Row 2 is the root "Extension" element; its Short/Definition text is
# updated to describe this specific extension instead of the generic
# "Extension" / "An Extension" placeholder text.
$elements.Range("K2").Value = "Episode Derived Code"
$elements.Range("L2").Value = "Derived code for the episode of care"
